$d = $word.ActiveDocument

$replacements = @(
    @("2024-01-26 Friday", "2024-01-27 Saturday"),
    @("312×2=", "544×8="),
    @("977×7=", "753×2="),
    @("732×9=", "768×7="),
    @("547×2=", "165×7="),
    @("661×6=", "713×3="),
    @("723×3=", "134×8="),
    @("967×5=", "364×5="),
    @("849×2=", "970×9="),
    @("128×8=", "413×3="),
    @("865×6=", "706×7="),
    @("882×4=", "609×9="),
    @("518×3=", "815×3="),
    @("359×5=", "465×3="),
    @("915×7=", "742×3="),
    @("442×3=", "888×7="),
    @("712×3=", "306×6="),
    @("813×9=", "111×4="),
    @("320×7=", "708×6="),
    @("853×9=", "472×3="),
    @("359×6=", "788×2="),
    @("986×7=", "175×3="),
    @("655×4=", "333×4="),
    @("517×5=", "585×8="),
    @("956×6=", "136×2="),
    @("748×8=", "329×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
